$wb = $excel.ActiveWorkbook

# Sheet "展览": 南宁·2024良牙动漫秋季盛典（秋典） 想去人数 1036 -> 1049
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1049

# Sheet "展览": 南宁·万圣漫控嘉年华10 想去人数 304 -> 305
$ws1.Range("F5").Value = 305

# Sheet "全部类型": 南宁·2024良牙动漫秋季盛典（秋典） 想去人数 1036 -> 1049
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1049

# Sheet "全部类型": 南宁·万圣漫控嘉年华10 想去人数 304 -> 305
$ws4.Range("F6").Value = 305
